# Added scalable phases function
# Update header row (row 1) of the active sheet: rename/reorder the
# existing header labels and append two new trailing headers so the
# row grows from columns A:AB to A:AD.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "time"
$ws.Range("B1").Value = "rew"
$ws.Range("C1").Value = "waitingTime"
$ws.Range("D1").Value = "Ep_num_train_rollouts"
$ws.Range("E1").Value = "Ep_rollout_length"
$ws.Range("F1").Value = "Ep_eval_freq"
$ws.Range("G1").Value = "Ep_eval_num_eps"
$ws.Range("H1").Value = "Ep_max_ep_steps"
$ws.Range("I1").Value = "Ep_generation_ep_steps"
$ws.Range("J1").Value = "Ep_test_num_eps"
$ws.Range("K1").Value = "A_agent_type"
$ws.Range("L1").Value = "A_single_agent"
$ws.Range("M1").Value = "P_gae_tau"
$ws.Range("N1").Value = "P_entropy_weight"
$ws.Range("O1").Value = "P_minibatch_size"
$ws.Range("P1").Value = "P_optimization_epochs"
$ws.Range("Q1").Value = "P_ppo_ratio_clip"
$ws.Range("R1").Value = "P_discount"
$ws.Range("S1").Value = "P_learning_rate"
$ws.Range("T1").Value = "P_clip_grads"
$ws.Range("U1").Value = "P_gradient_clip"
$ws.Range("V1").Value = "P_value_loss_coef"
$ws.Range("W1").Value = "R_rule_set"
$ws.Range("X1").Value = "R_rule_set_params"
$ws.Range("Y1").Value = "En_shape"
$ws.Range("Z1").Value = "En_rush_hour"
$ws.Range("AA1").Value = "En_uniform_generation_probability"
$ws.Range("AB1").Value = "M_reward_interpolation"
$ws.Range("AC1").Value = "M_state_interpolation"
$ws.Range("AD1").Value = "P_num_workers"
